$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0005255733983358369
$ws.Range("C2").Value = 0.0006956399435875938
$ws.Range("D2").Value = 0.0002342874786700122
$ws.Range("E2").Value = 0.0001725392614025623
$ws.Range("F2").Value = 0.0001745209208456799
$ws.Range("B3").Value = 0.0005465241411002353
$ws.Range("C3").Value = 0.0007245875801891088
$ws.Range("D3").Value = 0.0002455608404125087
$ws.Range("E3").Value = 0.0001735049794660881
$ws.Range("F3").Value = 0.000179678279964719
$ws.Range("B4").Value = 0.0006106475047999993
$ws.Range("C4").Value = 0.0008222499751718715
$ws.Range("D4").Value = 0.0002710141602437943
$ws.Range("E4").Value = 0.0001838725194102153
$ws.Range("F4").Value = 0.0001960533406236209
$ws.Range("B5").Value = 0.00061323837202508
$ws.Range("C5").Value = 0.0007727684028213844
$ws.Range("D5").Value = 0.0002616132993716746
$ws.Range("E5").Value = 0.0001832933802506886
$ws.Range("F5").Value = 0.0001839850610122084
$ws.Range("B6").Value = 0.0006715451017953455
$ws.Range("C6").Value = 0.0008236432867124677
$ws.Range("D6").Value = 0.0002725849786656909
$ws.Range("E6").Value = 0.0001856482806033455
$ws.Range("F6").Value = 0.0001917849620804191
$ws.Range("B7").Value = 0.001330756686511449
$ws.Range("C7").Value = 0.0008985550364013762
$ws.Range("D7").Value = 0.0002851666577043943
$ws.Range("E7").Value = 0.0001917608804069459
$ws.Range("F7").Value = 0.0002058649994432926
$ws.Range("B8").Value = 0.001419190007727593
$ws.Range("C8").Value = 0.001108905853470787
$ws.Range("D8").Value = 0.0002974100204301067
$ws.Range("E8").Value = 0.0001931234402582049
$ws.Range("F8").Value = 0.0002109799810568802
$ws.Range("B9").Value = 0.002052134192781523
$ws.Range("C9").Value = 0.001174770799116231
$ws.Range("D9").Value = 0.0002951232623308897
$ws.Range("E9").Value = 0.000191466678516008
$ws.Range("F9").Value = 0.0002127883594948799
$ws.Range("B10").Value = 0.002096100819180719
$ws.Range("C10").Value = 0.001161795758525841
$ws.Range("D10").Value = 0.0002915658225538209
$ws.Range("E10").Value = 0.0002108740992844105
$ws.Range("F10").Value = 0.0002049824403366074
$ws.Range("B11").Value = 0.0035353924613446
$ws.Range("C11").Value = 0.001484045895631425
$ws.Range("D11").Value = 0.0003117833018768579
$ws.Range("E11").Value = 0.0002032782990136184
$ws.Range("F11").Value = 0.0002127017205930315
$ws.Range("B12").Value = 0.007560601613367908
$ws.Range("C12").Value = 0.003513533374061808
$ws.Range("D12").Value = 0.0002853217202937231
$ws.Range("E12").Value = 0.0002001441601896659
$ws.Range("F12").Value = 0.0002207475615432486
$ws.Range("B13").Value = 0.01639937072061002
$ws.Range("C13").Value = 0.005661338327918201
$ws.Range("D13").Value = 0.0002906525204889476
$ws.Range("E13").Value = 0.0002057909205905162
$ws.Range("F13").Value = 0.0002181617222959176
$ws.Range("B14").Value = 0.02277728251414373
$ws.Range("C14").Value = 0.009765807497315109
$ws.Range("D14").Value = 0.0003348342009121552
$ws.Range("E14").Value = 0.0002080366795416921
$ws.Range("F14").Value = 0.0002442875411361456
$ws.Range("B15").Value = 0.03828838171204552
$ws.Range("C15").Value = 0.01963164651591796
$ws.Range("D15").Value = 0.0003139883774565533
$ws.Range("E15").Value = 0.0001707467384403571
$ws.Range("F15").Value = 0.0001789108195225708
$ws.Range("B16").Value = 0.07465762879641261
$ws.Range("C16").Value = 0.03163599542807788
$ws.Range("D16").Value = 0.000248636600736063
$ws.Range("E16").Value = 0.000139395019505173
$ws.Range("F16").Value = 0.0001567224599421024
$ws.Range("B17").Value = 0.1053145892790053
$ws.Range("C17").Value = 0.04680061661754735
$ws.Range("D17").Value = 0.0002727658793446608
$ws.Range("E17").Value = 0.0001478707807837054
$ws.Range("F17").Value = 0.0001438408193644136
$ws.Range("B18").Value = 0.1553878969780635
$ws.Range("C18").Value = 0.06007049338310026
$ws.Range("D18").Value = 0.0002565907582174987
$ws.Range("E18").Value = 0.0001459799610893242
$ws.Range("F18").Value = 0.000165651701099705
$ws.Range("B19").Value = 0.2137059942859924
$ws.Range("C19").Value = 0.08185323432320729
$ws.Range("D19").Value = 0.0002343049211776815
$ws.Range("E19").Value = 0.0001243582618189976
$ws.Range("F19").Value = 0.0001344591812812723
$ws.Range("B20").Value = 0.3651776934834197
$ws.Range("C20").Value = 0.1190055196243338
$ws.Range("D20").Value = 0.0002027933002682403
$ws.Range("E20").Value = 0.0001061716198455542
$ws.Range("F20").Value = 0.0001485299412161112
$ws.Range("B21").Value = 0.5197520854510367
$ws.Range("C21").Value = 0.1618393243907485
$ws.Range("D21").Value = 0.0002756824588868767
$ws.Range("E21").Value = 0.0001215233988477848
$ws.Range("F21").Value = 0.0001584417390404269
